$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Insert()
$ws.Range("A2:B2").Font.Bold = $false
$ws.Range("A2").Value = "Fix normals problem on cylinder"
$ws.Range("B2").Value = 3
$ws.Range("B3").Value = 4
